$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F (shifts the existing F column -- "LEITOSTOAL_UTI" --
# one place to the right, into G), and populate the freed-up F column with
# the new "PACIENTES_TOTAL" header + values.
$ws.Columns.Item(6).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

$ws.Range("F1").Value = "PACIENTES_TOTAL"

$pacientesTotal = @(30, 31, 31, 31, 31, 31, 37, 35, 33, 35, 40, 40, 40, 35, 39)
for ($i = 0; $i -lt $pacientesTotal.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $pacientesTotal[$i]
}

# New row 17 (previously unused) with the full set of values.
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 23
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 41
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 80
